$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @("30", "05", "0005", "Aegis Bloom", "ward", "BURST", "36", "4", "4.0", "20", "Deploys a bloom granting +60 shield for 8s.", "fx/relics/aegis_bloom.png", "sfx/relics/aegis.wav")

$targetRow = $ws.Range("A10:M10")
$targetRow.NumberFormat = "@"

$col = 1
foreach ($val in $rowData) {
    $ws.Cells.Item(10, $col).Value = $val
    $col++
}
